# Apply the "UI with relay settings" update to texts.xlsx
#
# Sheet "Typography": font file swapped from verdana.ttf to isocpeur.ttf for
# the "Default" and "Large" typographies, Large size changed 40 -> 30, and
# both rows get a new Wildcard Ranges value of "0-9".
#
# Sheet "Translation": a new "SI" (Slovenian) language column is appended,
# and 13 new rows of translated UI text (relay setup screens) are added.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Typography")
$ws2 = $wb.Worksheets.Item("Translation")

# ---- Typography sheet ----------------------------------------------------
$ws1.Range("C4").Value = "isocpeur.ttf"
$ws1.Range("I4").Value = "0-9"

$ws1.Range("C5").Value = "isocpeur.ttf"
$ws1.Range("D5").Value = 30
$ws1.Range("I5").Value = "0-9"

# ---- Translation sheet ----------------------------------------------------
# New header for the Slovenian ("SI") translation column
$ws2.Range("G3").Value = "SI"

# TEXT ID, TYPOGRAPHY NAME, ALIGNMENT, DIRECTION, GB, SI
$rows = @(
    @(4,  "SingleUseId3",  "Large",   "Center", "LTR", "Relay <value>",  "Rele <value>"),
    @(5,  "SingleUseId4",  "Large",   "Left",   "LTR", "0",              "0"),
    @(6,  "SingleUseId5",  "Default", "Center", "LTR", "<> ms",          "<> ms"),
    @(7,  "SingleUseId6",  "Default", "Center", "LTR", "<> ms",          "<> ms"),
    @(8,  "SingleUseId7",  "Default", "Left",   "LTR", "0",              "0"),
    @(9,  "SingleUseId8",  "Default", "Left",   "LTR", "0",              "0"),
    @(10, "SingleUseId10", "Default", "Center", "LTR", "Delay",          "Zamik"),
    @(11, "SingleUseId11", "Default", "Center", "LTR", "Duration",       "Dolzina"),
    @(12, "SingleUseId12", "Large",   "Left",   "LTR", "<digit>",        "<digit>"),
    @(13, "SingleUseId13", "Large",   "Left",   "LTR", "0",              "0"),
    @(14, "SingleUseId14", "Large",   "Left",   "LTR", "Relay 1 Setup",  "Rele 1"),
    @(15, "SingleUseId15", "Large",   "Left",   "LTR", "Relay 2 Setup",  "Rele 2"),
    @(16, "SingleUseId16", "Large",   "Left",   "LTR", "Relay 3 Setup",  "Rele 3")
)

# Cells whose translated text is the digit "0" need to stay TEXT (these are
# UI placeholder strings, not numeric values) -- assigning the literal "0"
# through .Value auto-coerces to a number, so those are written via a
# throwaway formula and converted back to a literal value/text afterwards.
function Set-TranslationCell($range, $text) {
    if ($text -eq "0") {
        $range.Formula = '="0"'
        $range.Copy() | Out-Null
        $range.PasteSpecial(-4163) | Out-Null   # xlPasteValues
    } else {
        $range.Value = $text
    }
}

foreach ($row in $rows) {
    $r = $row[0]
    Set-TranslationCell $ws2.Range("B$r") $row[1]
    Set-TranslationCell $ws2.Range("C$r") $row[2]
    Set-TranslationCell $ws2.Range("D$r") $row[3]
    Set-TranslationCell $ws2.Range("E$r") $row[4]
    Set-TranslationCell $ws2.Range("F$r") $row[5]
    Set-TranslationCell $ws2.Range("G$r") $row[6]
}
$excel.CutCopyMode = 0
